$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 7, pushing existing rows 7+ down by one.
$ws.Rows(7).Insert()

# Fill in the two previously-blank (but already hyperlink-styled) cells
# with plain-text URLs (no live hyperlink object).
$ws.Range("A5").Value = "https://en.wikipedia.org/wiki/Precession#:~:text=than%20this%2C%20however.-,Relativistic%20(Einsteinian),near%20a%20large%20rotating%20mass."
$ws.Range("A6").Value = "https://pressbooks.online.ucf.edu/osuniversityphysics/chapter/11-3-precession-of-a-gyroscope/#:~:text=The%20precessional%20angular%20velocity%20is%20given%20by%20%CF%89P%3Dr,frequency%20of%20the%20gyroscope%20disk."

# Populate the newly inserted row with a real, clickable hyperlink.
$ws.Hyperlinks.Add($ws.Range("A7"), "https://www.youtube.com/watch?v=ty9QSiVC2g0")
$ws.Range("A7").Style = "Hyperlink"

# Selection moves to A8 in the committed file.
$ws.Range("A8").Select() | Out-Null
